# Regenerated experiment task-order sheets.
# For each (fixed) tab position, the sheet is renamed and its task_order
# rows (column A = index, column B = stimulus/condition file name) are
# rewritten to a freshly generated order. Tab positions keep their
# physical slot (rId1..rId5) but the logical sheet now belongs to a
# different task.

$wb = $excel.ActiveWorkbook

function Set-TaskOrderSheet {
    param(
        $Worksheet,
        [string]$NewName,
        [object[]]$Rows
    )

    $Worksheet.Name = $NewName

    $oldLastRow = $Worksheet.UsedRange.Rows.Count
    if ($oldLastRow -lt 1) { $oldLastRow = 1 }
    $newLastRow = $Rows.Count + 1

    # Write header + data rows.
    $Worksheet.Cells.Item(1, 2).Value = "task_order"
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $i + 2
        $pair = $Rows[$i]
        $Worksheet.Cells.Item($r, 1).Value = $pair[0]
        $Worksheet.Cells.Item($r, 2).Value = $pair[1]
    }

    if ($newLastRow -gt $oldLastRow) {
        # Sheet grew: copy the index-column style (bold/centered/bordered)
        # from an already-styled cell down onto the newly added rows.
        $Worksheet.Range("A2").Copy()
        $fillRange = $Worksheet.Range("A" + ($oldLastRow + 1) + ":A" + $newLastRow)
        $fillRange.PasteSpecial(-4122)
    } elseif ($newLastRow -lt $oldLastRow) {
        # Sheet shrank: remove the now-unused trailing rows entirely so the
        # used range / dimension shrinks too.
        $deleteRange = $Worksheet.Rows(($newLastRow + 1).ToString() + ":" + $oldLastRow.ToString())
        $deleteRange.Delete()
    }
}

Set-TaskOrderSheet $wb.Worksheets.Item(1) "vSAT_TO-16515890462110953" @(
    , @(0, "SAT_stims-16515890461642554.csv")
    , @(1, "vSAT_stims-165158904619547.csv")
    , @(2, "SAT_stims-1651589046148629.csv")
    , @(3, "vSAT_stims-16515890461798472.csv")
)

Set-TaskOrderSheet $wb.Worksheets.Item(2) "NB_TO-1651589047682658" @(
    , @(0, "ZB-match_2-16515890462579937.csv")
    , @(1, "TB-16515890473522236.csv")
    , @(2, "ZB-match_7-1651589046413727.csv")
    , @(3, "TB-16515890473835163.csv")
    , @(4, "OB-1651589046444978.csv")
    , @(5, "TB-16515890476670325.csv")
    , @(6, "OB-16515890468221488.csv")
    , @(7, "ZB-match_9-1651589046226721.csv")
    , @(8, "OB-16515890469961047.csv")
)

Set-TaskOrderSheet $wb.Worksheets.Item(3) "TOL_TO-1651589047729506" @(
    , @(0, "MM_stims-1651589047698256.csv")
    , @(1, "ZM_stims-1651589047682658.csv")
    , @(2, "MM_stims-16515890477139227.csv")
    , @(3, "ZM_stims-1651589047698256.csv")
    , @(4, "MM_stims-1651589047729506.csv")
    , @(5, "ZM_stims-16515890477139227.csv")
)

Set-TaskOrderSheet $wb.Worksheets.Item(4) "GNG_TO-16515890477607546" @(
    , @(0, "go_stims-1651589047729506.csv")
    , @(1, "GNG_stims-165158904774513.csv")
    , @(2, "go_stims-165158904774513.csv")
    , @(3, "GNG_stims-16515890477607546.csv")
)

Set-TaskOrderSheet $wb.Worksheets.Item(5) "RS_TO-16515890477607546" @(
    , @(0, "eyes open")
    , @(1, "eyes closed")
)
